# Filtered movement data to match siminputrow matrix
# Rows 39-57 were labeled "Feb" (id_month) but should be "Jan"
# Rows 58-65 were labeled "Feb" (id_month) but should be "Mar"
# (No other rows' id_month text actually changes; the "Feb" shared string
#  becomes unused and is dropped on save, shifting later string indices.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 39; $r -le 57; $r++) {
    $ws.Range("C$r").Value = "Jan"
}

for ($r = 58; $r -le 65; $r++) {
    $ws.Range("C$r").Value = "Mar"
}
